$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.259.28'
$ws.Range('E2').Value = '  -7.10%  '
$ws.Range('D3').Value = '3.296.37'
$ws.Range('E3').Value = '  -4.05%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''555.67'
$ws.Range('E5').Value = '  -4.41%  '
$ws.Range('D6').Value = '''127.30'
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.292.94'
$ws.Range('E8').Value = '  -4.11%  '
$ws.Range('D9').Value = '''0.467'
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('E10').Value = '  -3.93%  '
$ws.Range('D11').Value = '''0.116'
$ws.Range('E11').Value = '  -5.42%  '
$ws.Range('D12').Value = '''0.369'
$ws.Range('E12').Value = '  -3.46%  '
$ws.Range('D13').Value = '3.857.80'
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').Value = '3.289.98'
$ws.Range('E15').Value = '  -4.45%  '
$ws.Range('E16').Value = '  -5.87%  '
$ws.Range('D17').Value = '''23.99'
$ws.Range('E17').Value = '  -3.34%  '
$ws.Range('D18').Value = '59.447.62'
$ws.Range('E18').Value = '  -6.67%  '
$ws.Range('D19').Value = '''5.61'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '''13.16'
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('D21').Value = '''8.87'
$ws.Range('E21').Value = '  -9.88%  '
$ws.Range('D22').Value = '''349.51'
$ws.Range('E22').Value = '  -8.92%  '
$ws.Range('D23').Value = '''0.551'
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = '3.427.88'
$ws.Range('E25').Value = '  -4.25%  '
$ws.Range('D26').Value = '''68.24'
$ws.Range('E26').Value = '  -7.79%  '
$ws.Range('D27').Value = '''0.0000109'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').Value = '''0.998'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '''7.22'
$ws.Range('E29').Value = '  +2.79%  '
$ws.Range('D30').Value = '''1.46'
$ws.Range('E30').Value = '  +2.88%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = '''0.151'
$ws.Range('E31').Value = '  -2.19%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''7.75'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '''2.08'
$ws.Range('E33').Value = '  -5.83%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '''1.00'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '3.327.43'
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('D36').Value = '''22.60'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').Value = '''5.30'
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').Value = '''157.58'
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('D41').Value = '''0.0742'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('D42').Value = '''0.998'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').Value = '''40.60'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').Value = '''1.17'
$ws.Range('E44').Value = '  +5.20%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '''4.26'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '''0.737'
$ws.Range('E46').Value = '  -6.82%  '
$ws.Range('D47').Value = '''22.54'
$ws.Range('E47').Value = '  -3.95%  '
$ws.Range('E48').Value = '  -4.77%  '
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').Value = '''2.36'
$ws.Range('E50').Value = '  +14.50%  '
$ws.Range('D51').Value = '''21.60'
$ws.Range('E51').Value = '  +6.59%  '
